$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: set a cell's value as TEXT (shared string) while preserving the
# cell's original number format (so cells whose style keeps its original
# numFmtId, e.g. "1" in a numeric-formatted column, are not silently
# converted to a number by the auto-detection logic).
# ---------------------------------------------------------------------
function Set-TextValue {
    param($cell, [string]$text)
    $origFormat = $cell.NumberFormat
    $cell.NumberFormat = "@"
    $cell.Value = $text
    if ($origFormat -ne "General") {
        $cell.NumberFormat = $origFormat
    }
}

# ---------------------------------------------------------------------
# A new row of data (row 9) is inserted before the old footer row, which
# shifts the old row 9 (footer: date/page/credit) down to row 10.
# ---------------------------------------------------------------------
$ws.Rows.Item(9).Insert()
$ws.Rows.Item(9).RowHeight = 25.5

# The old row 8 only contained a merged "P8:Q8" pair; unmerge it now so
# every cell in row 8 can be addressed/written individually.
$ws.Range("P8:Q8").UnMerge()

# Seed row 8 (A:Q) with row 7's cell formatting, since row 8 previously
# had no cells at all in columns A-O.
$ws.Range("A7:Q7").Copy()
$ws.Range("A8:Q8").PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------
# Row 7 : first sold item line
# ---------------------------------------------------------------------
$ws.Cells.Item(7,1).Value = 1                                  # A7
Set-TextValue $ws.Cells.Item(7,3)  "BRONCHICUM ELIXIR 100 ML"   # C7
Set-TextValue $ws.Cells.Item(7,8)  "1:0"                        # H7
Set-TextValue $ws.Cells.Item(7,12) "1"                          # L7
Set-TextValue $ws.Cells.Item(7,14) "69.00"                      # N7
Set-TextValue $ws.Cells.Item(7,16) "69.0000"                    # P7
Set-TextValue $ws.Cells.Item(7,17) "1:0"                        # Q7

# ---------------------------------------------------------------------
# Row 8 : second sold item line
# ---------------------------------------------------------------------
$ws.Cells.Item(8,1).Value = 2                                                # A8
Set-TextValue $ws.Cells.Item(8,3)  "CONTAFEVER N 200MG/5ML SUSP. 120ML"      # C8
Set-TextValue $ws.Cells.Item(8,8)  "11:0"                                    # H8
Set-TextValue $ws.Cells.Item(8,12) "1"                                       # L8
Set-TextValue $ws.Cells.Item(8,14) "33.00"                                   # N8
Set-TextValue $ws.Cells.Item(8,16) "33.0000"                                 # P8
Set-TextValue $ws.Cells.Item(8,17) "1:0"                                     # Q8

# ---------------------------------------------------------------------
# Row 9 (new) : transaction-count total
# ---------------------------------------------------------------------
$ws.Cells.Item(9,16).Value = 102   # P9

# ---------------------------------------------------------------------
# Merge cells: row 7 layout (A:B, C:G, H:K, L:M, N:O) now also applies to
# row 8, and the new total row gets its own "P9:Q9" merge.
# ---------------------------------------------------------------------
$ws.Range("A8:B8").Merge() | Out-Null
$ws.Range("C8:G8").Merge() | Out-Null
$ws.Range("H8:K8").Merge() | Out-Null
$ws.Range("L8:M8").Merge() | Out-Null
$ws.Range("N8:O8").Merge() | Out-Null
$ws.Range("P9:Q9").Merge() | Out-Null
